# Change the presentation's theme font scheme from "Calibri Light-Constantia"
# to the built-in "Cambria" font scheme (major/minor Latin typeface -> Cambria).
$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$theme = $m.Theme
$fs = $theme.ThemeFontScheme

# Font scheme / theme-fonts display name (best effort - PowerPoint's object
# model normally derives this name, but set it explicitly in case it sticks).
$fs.Name = "Cambria"

# Major (heading) font: Calibri Light -> Cambria
$major = $fs.MajorFont
$major.Latin = "Cambria"

# Minor (body) font: Constantia -> Cambria
$minor = $fs.MinorFont
$minor.Latin = "Cambria"
